# "Sesiones en otra hoja de excel"
# Move the "Sesion" mini-table (previously columns N:P on Hoja1) onto a
# brand-new "Hoja2", and extend it with two more session rows
# (Informatica 4 / Informatica 5). Hoja2 becomes the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed right after Hoja1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

# Cut the existing Sesion table (N1:P3) off Hoja1 and paste it at A1 on
# Hoja2 - this carries the values, shared-string usage and cell styles
# along in one step.
$ws1.Range("N1:P3").Cut($ws2.Range("A1"))

# Cut() leaves the (now empty) source cells behind with their old
# formatting; clear that out so Hoja1's used range shrinks back down to
# A1:M3, same as the rest of the sheet.
$ws1.Range("N1:P3").Clear()

# Hoja1 column P had a custom "best fit" width; carry an equivalent
# custom width over onto Hoja2 column C, which now holds the same data.
$ws2.Columns.Item(3).ColumnWidth = 11.67

# Two new session rows for Informatica 4 and Informatica 5.
$ws2.Range("A4").Value = "IS"
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = "Informatica 4"

$ws2.Range("A5").Value = "IS"
$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = "Informatica 5"

# Match the selections/active sheet left behind by the edit.
$ws1.Range("N1:P3").Select() | Out-Null

$ws2.Range("F3").Select() | Out-Null
$ws2.Activate() | Out-Null
